$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -4
$ws.Range("F6").Value = -6
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = 11
$ws.Range("F13").Value = -2
$ws.Range("F14").Value = -7
$ws.Range("F16").Value = -1
